# Extend the "In Class Demonstration" regression-summary table from 3
# dependent-variable columns (LF, C/A, FFR) out to 6 (adding LF Lag,
# C/A Lag, FFR Lag), and refresh the coefficient values that were already
# present for the first three columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 : header labels -- copy the bold/centered/bordered header
#      style already used by B1:D1 onto the new header cells first ----
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:G1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("E1").Value = "LF Lag"
$ws.Range("F1").Value = "C/A Lag"
$ws.Range("G1").Value = "FFR Lag"

# ---- Row 2 : LF Lag coefficients ----
# (values that are plain decimals need a leading apostrophe so they are
# entered as text, exactly like the pre-existing "0.043"/"0.029" cells,
# instead of being auto-recognised as numbers)
$ws.Range("C2").Value = "'0.486"
$ws.Range("D2").Value = "'0.217"
$ws.Range("E2").Value = "1.0***"
$ws.Range("F2").Value = "-0.0**"
$ws.Range("G2").Value = "'0.0"

# ---- Row 3 : C/A Lag coefficients ----
$ws.Range("C3").Value = "-0.478*"
$ws.Range("D3").Value = "'-0.068"
$ws.Range("E3").Value = "-0.0**"
$ws.Range("F3").Value = "1.0***"
$ws.Range("G3").Value = "-0.0***"

# ---- Row 4 : FFR Lag coefficients ----
$ws.Range("C4").Value = "0.742***"
$ws.Range("D4").Value = "0.893***"
$ws.Range("E4").Value = "'-0.0"
$ws.Range("F4").Value = "-0.0*"
$ws.Range("G4").Value = "1.0***"

# ---- Row 5 : Constant ----
$ws.Range("C5").Value = "'-0.535"
$ws.Range("D5").Value = "-1.493***"
$ws.Range("E5").Value = "'0.0"
$ws.Range("F5").Value = "'-0.0"
$ws.Range("G5").Value = "0.0*"

# ---- Row 6 : r2_adj (plain numbers) ----
$ws.Range("C6").Value = 0.4
$ws.Range("D6").Value = 0.91
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
